# Update the "想去人数" (want-to-go count) and "最低票价" (min price) numbers
# in the 北京-漫展信息 workbook to the values regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$ws1.Range("F6").Value = 271
$ws1.Range("F7").Value = 13136
$ws1.Range("G7").Value = 110
$ws1.Range("F8").Value = 65
$ws1.Range("F10").Value = 280
$ws1.Range("F11").Value = 3990
$ws1.Range("F12").Value = 6688
$ws1.Range("F15").Value = 3507
$ws1.Range("F16").Value = 40
$ws1.Range("F17").Value = 170
$ws1.Range("F21").Value = 128
$ws1.Range("F22").Value = 3654
$ws1.Range("F25").Value = 3615
$ws1.Range("F26").Value = 3615
$ws1.Range("F28").Value = 1917
$ws1.Range("F29").Value = 107
$ws1.Range("F30").Value = 241
$ws1.Range("F31").Value = 6842
$ws1.Range("F34").Value = 1690
$ws1.Range("F35").Value = 2028
$ws1.Range("F37").Value = 111
$ws1.Range("F38").Value = 1080
$ws1.Range("F40").Value = 222
$ws1.Range("F41").Value = 13
$ws1.Range("F45").Value = 4
$ws1.Range("F47").Value = 1226
$ws1.Range("F48").Value = 1825
$ws1.Range("F49").Value = 71
$ws1.Range("F50").Value = 165

# --- Sheet "本地生活" ---
$ws3.Range("F2").Value = 465
$ws3.Range("F3").Value = 631

# --- Sheet "全部类型" ---
$ws4.Range("F6").Value = 465
$ws4.Range("F7").Value = 631
$ws4.Range("F9").Value = 271
$ws4.Range("F10").Value = 13136
$ws4.Range("G10").Value = 110
$ws4.Range("F11").Value = 65
$ws4.Range("F14").Value = 280
$ws4.Range("F15").Value = 3991
$ws4.Range("F16").Value = 6688
$ws4.Range("F18").Value = 3507
$ws4.Range("F19").Value = 40
$ws4.Range("F25").Value = 128
$ws4.Range("F28").Value = 3615
$ws4.Range("F30").Value = 107
$ws4.Range("F31").Value = 241
$ws4.Range("F32").Value = 6843
$ws4.Range("F36").Value = 1690
$ws4.Range("F37").Value = 2028
$ws4.Range("F39").Value = 111
$ws4.Range("F40").Value = 1080
$ws4.Range("F41").Value = 222
$ws4.Range("F47").Value = 1825
$ws4.Range("F48").Value = 71
$ws4.Range("F50").Value = 165
